# Applies the "new inventory data" edit described by the commit:
#  1. Remove the empty "Sheet2" worksheet (Sheet3 becomes the 2nd sheet).
#  2. Add new inventory metadata columns (M:X) on row 1 (headers) and
#     populate the first data row (row 2) with sample inventory values.
#  3. Widen the existing columns A:K (and the generic overflow column)
#     to fit the new, longer header/content text.
#  4. Update the active selection to S2 and keep gridlines visible.

$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inventory")

# --- 1. Delete the unused "Sheet2" tab -------------------------------------
$wb.Worksheets.Item("Sheet2").Delete()

# --- 2. New header cells (row 1, columns M:X) ------------------------------
$ws.Range("M1").Value = "InventoryType"
$ws.Range("N1").Value = "Category"
$ws.Range("O1").Value = "Department"
$ws.Range("P1").Value = "Brand"
$ws.Range("Q1").Value = "Item Name"
$ws.Range("R1").Value = "Short Name"
$ws.Range("S1").Value = "Item Type"
$ws.Range("T1").Value = "UOM Purchase"
$ws.Range("U1").Value = "Sales UOM"
$ws.Range("V1").Value = "Stock Type"
$ws.Range("W1").Value = "Vendor"
$ws.Range("X1").Value = "Floor"

# --- New inventory record (row 2, columns M:S and W) -----------------------
$ws.Range("M2").Value = "inventory"
$ws.Range("N2").Value = "BISCUITS"
$ws.Range("O2").Value = "SNACKS"
$ws.Range("P2").Value = "PARLE"
$ws.Range("Q2").Value = "BISC Test Parle"
$ws.Range("R2").Value = "PARLEs"
$ws.Range("S2").Value = "BULK"
$ws.Range("W2").Value = "AKSHAYA TRADERS"

# --- 3. Resize columns to match the new (wider) content ---------------------
$ws.Columns.Item(1).ColumnWidth = 12.166666666666666
$ws.Columns.Item(2).ColumnWidth = 12.166666666666666
$ws.Columns.Item(3).ColumnWidth = 19.5
$ws.Columns.Item(4).ColumnWidth = 26.5
$ws.Columns.Item(5).ColumnWidth = 23.833333333333336
$ws.Columns.Item(6).ColumnWidth = 17.333333333333336
$ws.Columns.Item(7).ColumnWidth = 14.666666666666666
$ws.Columns.Item(8).ColumnWidth = 17.833333333333336
$ws.Columns.Item(9).ColumnWidth = 14.666666666666666
$ws.Columns.Item(10).ColumnWidth = 25.833333333333336
$ws.Columns.Item(11).ColumnWidth = 19.666666666666668
$ws.Range("L1:X1048576").ColumnWidth = 11.5

# --- 4. View: keep gridlines on, select S2, scroll toward the new columns ---
$excel.ActiveWindow.DisplayGridlines = $true
$ws.Range("S2").Select()
